$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts N->O, O->P, P->Q)
$ws.Columns("N:N").Insert()

# Make "Repayment Schedule" the active/selected sheet and set its selection
$ws.Activate() | Out-Null
$ws.Range("L17").Select() | Out-Null
